# Scheduled runner update: refresh computed market-price / leve-profit
# columns (currentAveragePrice*, LevePrice*, LeveProfit*) for the rows
# whose underlying Universalis price data changed.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 164.85715
$ws.Range("I39").Value = 44.384617
$ws.Range("J39").Value = 360.625
$ws.Range("K39").Value = 133.153851
$ws.Range("L39").Value = 1081.875
$ws.Range("M39").Value = 162.846149
$ws.Range("N39").Value = -1673.875
$ws.Range("H86").Value = 6928.5293
$ws.Range("I86").Value = 12997.5
$ws.Range("J86").Value = 5061.154
$ws.Range("K86").Value = 12997.5
$ws.Range("L86").Value = 5061.154
$ws.Range("M86").Value = -11874.5
$ws.Range("N86").Value = -7307.154
$ws.Range("H89").Value = 6928.5293
$ws.Range("I89").Value = 12997.5
$ws.Range("J89").Value = 5061.154
$ws.Range("K89").Value = 64987.5
$ws.Range("L89").Value = 25305.77
$ws.Range("M89").Value = -59371.5
$ws.Range("N89").Value = -36537.77
$ws.Range("H121").Value = 2724.5
$ws.Range("J121").Value = 2724.5
$ws.Range("L121").Value = 8173.5
$ws.Range("N121").Value = -11667.5
$ws.Range("H137").Value = 53713.688
$ws.Range("I137").Value = 73333.32000000001
$ws.Range("K137").Value = 219999.96
$ws.Range("M137").Value = -217449.96
$ws.Range("H138").Value = 3325.3076
$ws.Range("I138").Value = 2152.8125
$ws.Range("J138").Value = 3708.1633
$ws.Range("K138").Value = 6458.4375
$ws.Range("L138").Value = 11124.4899
$ws.Range("M138").Value = -1318.4375
$ws.Range("N138").Value = -21404.4899

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8166.216
$ws.Range("I32").Value = 6733.0127
$ws.Range("J32").Value = 22498.25
$ws.Range("K32").Value = 6733.0127
$ws.Range("L32").Value = 22498.25
$ws.Range("M32").Value = -6446.0127
$ws.Range("N32").Value = -23072.25
$ws.Range("H45").Value = 7940130.5
$ws.Range("I45").Value = 15874096
$ws.Range("K45").Value = 15874096
$ws.Range("M45").Value = -15873719
$ws.Range("H46").Value = 3082.2856
$ws.Range("I46").Value = 3141.6667
$ws.Range("J46").Value = 3037.75
$ws.Range("K46").Value = 3141.6667
$ws.Range("L46").Value = 3037.75
$ws.Range("M46").Value = -2822.6667
$ws.Range("N46").Value = -3675.75
$ws.Range("H122").Value = 2537044
$ws.Range("I122").Value = 2926013.8
$ws.Range("K122").Value = 8778041.399999999
$ws.Range("M122").Value = -8775591.399999999
$ws.Range("H132").Value = 37893.52
$ws.Range("I132").Value = 1685.6
$ws.Range("J132").Value = 182725.2
$ws.Range("K132").Value = 5056.799999999999
$ws.Range("L132").Value = 548175.6000000001
$ws.Range("M132").Value = -2526.799999999999
$ws.Range("N132").Value = -553235.6000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24946.904
$ws.Range("I31").Value = 2742.8928
$ws.Range("J31").Value = 69354.92999999999
$ws.Range("K31").Value = 2742.8928
$ws.Range("L31").Value = 69354.92999999999
$ws.Range("M31").Value = -2447.8928
$ws.Range("N31").Value = -69944.92999999999
$ws.Range("H34").Value = 24946.904
$ws.Range("I34").Value = 2742.8928
$ws.Range("J34").Value = 69354.92999999999
$ws.Range("K34").Value = 2742.8928
$ws.Range("L34").Value = 69354.92999999999
$ws.Range("M34").Value = -2540.8928
$ws.Range("N34").Value = -69758.92999999999
$ws.Range("H100").Value = 56243.75
$ws.Range("J100").Value = 56243.75
$ws.Range("L100").Value = 56243.75
$ws.Range("N100").Value = -58407.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 11000
$ws.Range("J11").Value = 10000
$ws.Range("L11").Value = 30000
$ws.Range("N11").Value = -30280
$ws.Range("H34").Value = 1315
$ws.Range("J34").Value = 2000
$ws.Range("L34").Value = 6000
$ws.Range("N34").Value = -6168
$ws.Range("H37").Value = 47740
$ws.Range("J37").Value = 47740
$ws.Range("L37").Value = 143220
$ws.Range("N37").Value = -143444
$ws.Range("H86").Value = 170
$ws.Range("I86").Value = 162.5
$ws.Range("J86").Value = 200
$ws.Range("K86").Value = 487.5
$ws.Range("L86").Value = 600
$ws.Range("M86").Value = 698.5
$ws.Range("N86").Value = -2972
$ws.Range("H89").Value = 170
$ws.Range("I89").Value = 162.5
$ws.Range("J89").Value = 200
$ws.Range("K89").Value = 1462.5
$ws.Range("L89").Value = 1800
$ws.Range("M89").Value = 4465.5
$ws.Range("N89").Value = -13656
$ws.Range("H107").Value = 496.85715
$ws.Range("I107").Value = 359.25
$ws.Range("J107").Value = 551.9
$ws.Range("K107").Value = 1077.75
$ws.Range("L107").Value = 1655.7
$ws.Range("M107").Value = 842.25
$ws.Range("N107").Value = -5495.7
$ws.Range("H121").Value = 1437.3889
$ws.Range("J121").Value = 1600.8125
$ws.Range("L121").Value = 4802.4375
$ws.Range("N121").Value = -7422.4375
$ws.Range("H136").Value = 6195.615
$ws.Range("I136").Value = 6054.3
$ws.Range("J136").Value = 6666.6665
$ws.Range("K136").Value = 18162.9
$ws.Range("L136").Value = 19999.9995
$ws.Range("M136").Value = -13062.9
$ws.Range("N136").Value = -30199.9995

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 38450.355
$ws.Range("I45").Value = 28000
$ws.Range("J45").Value = 44256.11
$ws.Range("K45").Value = 28000
$ws.Range("L45").Value = 44256.11
$ws.Range("M45").Value = -27441
$ws.Range("N45").Value = -45374.11
$ws.Range("H95").Value = 18686
$ws.Range("J95").Value = 18686
$ws.Range("L95").Value = 18686
$ws.Range("N95").Value = -24178
$ws.Range("H102").Value = 4354303.5
$ws.Range("I102").Value = 7937616
$ws.Range("J102").Value = 1403340.2
$ws.Range("K102").Value = 7937616
$ws.Range("L102").Value = 1403340.2
$ws.Range("M102").Value = -7935994
$ws.Range("N102").Value = -1406584.2
$ws.Range("H132").Value = 3130.22
$ws.Range("I132").Value = 2801.0256
$ws.Range("J132").Value = 4297.364
$ws.Range("K132").Value = 8403.076799999999
$ws.Range("L132").Value = 12892.092
$ws.Range("M132").Value = -5873.076799999999
$ws.Range("N132").Value = -17952.092

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6560.657
$ws.Range("I7").Value = 5883.05
$ws.Range("J7").Value = 7464.1333
$ws.Range("K7").Value = 5883.05
$ws.Range("L7").Value = 7464.1333
$ws.Range("M7").Value = -5771.05
$ws.Range("N7").Value = -7688.1333
$ws.Range("H40").Value = 7634.6665
$ws.Range("I40").Value = 3956.7
$ws.Range("K40").Value = 3956.7
$ws.Range("M40").Value = -3820.7
$ws.Range("H126").Value = 6560.657
$ws.Range("I126").Value = 5883.05
$ws.Range("J126").Value = 7464.1333
$ws.Range("K126").Value = 17649.15
$ws.Range("L126").Value = 22392.3999
$ws.Range("M126").Value = -15179.15
$ws.Range("N126").Value = -27332.3999
$ws.Range("H132").Value = 8738.6875
$ws.Range("I132").Value = 9205.286
$ws.Range("K132").Value = 27615.858
$ws.Range("M132").Value = -25085.858

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 12482571
$ws.Range("I132").Value = 13701308
$ws.Range("J132").Value = 1361604.5
$ws.Range("K132").Value = 41103924
$ws.Range("L132").Value = 4084813.5
$ws.Range("M132").Value = -41101394
$ws.Range("N132").Value = -4089873.5
